# Adjust stochastic example so Scenario A & B are actually different.
# Scenario A (sheet1) keeps its original representative-period weights;
# Scenario B (sheet2) gets a distinct set of weights for rows 8-14 (rp01-rp07).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioB")

$ws.Range("D8").Value  = 40
$ws.Range("D9").Value  = 125
$ws.Range("D10").Value = 51
$ws.Range("D11").Value = 48
$ws.Range("D12").Value = 43
$ws.Range("D13").Value = 20
$ws.Range("D14").Value = 37
